$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 234.6
$ws.Range("I9").Value = 295
$ws.Range("J9").Value = 144
$ws.Range("K9").Value = 295
$ws.Range("L9").Value = 144
$ws.Range("M9").Value = -126
$ws.Range("N9").Value = -482
$ws.Range("H20").Value = 2955.3333
$ws.Range("I20").Value = 395
$ws.Range("J20").Value = 15757
$ws.Range("K20").Value = 395
$ws.Range("L20").Value = 15757
$ws.Range("M20").Value = -165
$ws.Range("N20").Value = -16217
$ws.Range("H33").Value = 146.9
$ws.Range("I33").Value = 103.166664
$ws.Range("K33").Value = 103.166664
$ws.Range("M33").Value = 125.833336
$ws.Range("H35").Value = 2955.3333
$ws.Range("I35").Value = 395
$ws.Range("J35").Value = 15757
$ws.Range("K35").Value = 395
$ws.Range("L35").Value = 15757
$ws.Range("M35").Value = -16
$ws.Range("N35").Value = -16515
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 1808.7368
$ws.Range("I100").Value = 1462
$ws.Range("J100").Value = 2779.6
$ws.Range("K100").Value = 1462
$ws.Range("L100").Value = 2779.6
$ws.Range("M100").Value = -921
$ws.Range("N100").Value = -3861.6
$ws.Range("H104").Value = 766.2
$ws.Range("I104").Value = 766.2
$ws.Range("K104").Value = 2298.6
$ws.Range("M104").Value = -551.6000000000004
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H137").Value = 1912.8
$ws.Range("I137").Value = 1031.909
$ws.Range("J137").Value = 2989.4443
$ws.Range("K137").Value = 3095.727
$ws.Range("L137").Value = 8968.332900000001
$ws.Range("M137").Value = -545.7270000000003
$ws.Range("N137").Value = -14068.3329

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4169309.2
$ws.Range("I32").Value = 1026.3
$ws.Range("K32").Value = 1026.3
$ws.Range("M32").Value = -739.3
$ws.Range("H97").Value = 887.64703
$ws.Range("I97").Value = 660.8
$ws.Range("J97").Value = 1211.7142
$ws.Range("K97").Value = 660.8
$ws.Range("L97").Value = 1211.7142
$ws.Range("M97").Value = -164.8
$ws.Range("N97").Value = -2203.7142
$ws.Range("H110").Value = 2608.3
$ws.Range("J110").Value = 4135
$ws.Range("L110").Value = 4135
$ws.Range("N110").Value = -8225
$ws.Range("H132").Value = 1787.25
$ws.Range("I132").Value = 1803.3903
$ws.Range("K132").Value = 5410.1709
$ws.Range("M132").Value = -2880.1709
$ws.Range("H139").Value = 37715
$ws.Range("J139").Value = 37715
$ws.Range("L139").Value = 37715
$ws.Range("N139").Value = -47995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3600.8572
$ws.Range("I20").Value = 3600.8572
$ws.Range("K20").Value = 3600.8572
$ws.Range("M20").Value = -3353.8572
$ws.Range("H80").Value = 269.33334
$ws.Range("J80").Value = 325
$ws.Range("L80").Value = 325
$ws.Range("N80").Value = -2321
$ws.Range("H83").Value = 269.33334
$ws.Range("J83").Value = 325
$ws.Range("L83").Value = 1625
$ws.Range("N83").Value = -11609
$ws.Range("H86").Value = 4014.7693
$ws.Range("I86").Value = 1536.25
$ws.Range("K86").Value = 1536.25
$ws.Range("M86").Value = -413.25
$ws.Range("H89").Value = 4014.7693
$ws.Range("I89").Value = 1536.25
$ws.Range("K89").Value = 7681.25
$ws.Range("M89").Value = -2065.25
$ws.Range("H94").Value = 383.5
$ws.Range("I94").Value = 334.26666
$ws.Range("K94").Value = 334.26666
$ws.Range("M94").Value = 116.73334
$ws.Range("H107").Value = 5478.357
$ws.Range("I107").Value = 956.7143
$ws.Range("K107").Value = 956.7143
$ws.Range("M107").Value = 963.2857
$ws.Range("H134").Value = 1402.08
$ws.Range("I134").Value = 1223.0416
$ws.Range("K134").Value = 3669.1248
$ws.Range("M134").Value = -1134.1248

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8450
$ws.Range("I6").Value = 358
$ws.Range("K6").Value = 358
$ws.Range("M6").Value = -245
$ws.Range("H58").Value = 1361
$ws.Range("I58").Value = 1361
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1361
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1158
$ws.Range("N58").ClearContents()
$ws.Range("H80").Value = 47870.2
$ws.Range("J80").Value = 47870.2
$ws.Range("L80").Value = 47870.2
$ws.Range("N80").Value = -50116.2
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value = 47870.2
$ws.Range("J83").Value = 47870.2
$ws.Range("L83").Value = 143610.6
$ws.Range("N83").Value = -154842.6
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H136").Value = 1361
$ws.Range("I136").Value = 1361
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4083
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1533
$ws.Range("N136").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 16.916666
$ws.Range("I7").Value = 16.916666
$ws.Range("K7").Value = 50.749998
$ws.Range("M7").Value = 61.250002
$ws.Range("H80").Value = 3943.4
$ws.Range("I80").Value = 3798.6667
$ws.Range("J80").Value = 4160.5
$ws.Range("K80").Value = 11396.0001
$ws.Range("L80").Value = 12481.5
$ws.Range("M80").Value = -10460.0001
$ws.Range("N80").Value = -14353.5
$ws.Range("H83").Value = 3943.4
$ws.Range("I83").Value = 3798.6667
$ws.Range("J83").Value = 4160.5
$ws.Range("K83").Value = 34188.0003
$ws.Range("L83").Value = 37444.5
$ws.Range("M83").Value = -29508.0003
$ws.Range("N83").Value = -46804.5
$ws.Range("H140").Value = 2329.9
$ws.Range("I140").Value = 2005.4445
$ws.Range("K140").Value = 6016.333500000001
$ws.Range("M140").Value = -836.3335000000006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 99
$ws.Range("J2").Value = 36.5
$ws.Range("K2").Value = 99
$ws.Range("L2").Value = 36.5
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = -262.5
$ws.Range("H26").Value = 34999
$ws.Range("J26").Value = 34999
$ws.Range("L26").Value = 34999
$ws.Range("N26").Value = -35559
$ws.Range("H50").Value = 34999
$ws.Range("J50").Value = 34999
$ws.Range("L50").Value = 34999
$ws.Range("N50").Value = -35995
$ws.Range("H70").Value = 1067.1428
$ws.Range("I70").Value = 1123
$ws.Range("J70").Value = 992.6667
$ws.Range("K70").Value = 1123
$ws.Range("L70").Value = 992.6667
$ws.Range("M70").Value = -853
$ws.Range("N70").Value = -1532.6667
$ws.Range("H73").Value = 1067.1428
$ws.Range("I73").Value = 1123
$ws.Range("J73").Value = 992.6667
$ws.Range("K73").Value = 1123
$ws.Range("L73").Value = 992.6667
$ws.Range("M73").Value = -187
$ws.Range("N73").Value = -2864.6667
$ws.Range("H122").Value = 1960.6
$ws.Range("I122").Value = 1950.8334
$ws.Range("K122").Value = 5852.5002
$ws.Range("M122").Value = -3402.5002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 4275
$ws.Range("I19").Value = 4366.6665
$ws.Range("K19").Value = 4366.6665
$ws.Range("M19").Value = -4196.6665
$ws.Range("H46").Value = 2636
$ws.Range("I46").Value = 1522.5
$ws.Range("J46").Value = 3229.8667
$ws.Range("K46").Value = 1522.5
$ws.Range("L46").Value = 3229.8667
$ws.Range("M46").Value = -1334.5
$ws.Range("N46").Value = -3605.8667
$ws.Range("H136").Value = 2299.6667
$ws.Range("I136").Value = 2299.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6899.000100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4349.000100000001
$ws.Range("N136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 991.44446
$ws.Range("J96").Value = 1047
$ws.Range("L96").Value = 1047
$ws.Range("N96").Value = -3793
$ws.Range("H126").Value = 3249.6296
$ws.Range("I126").Value = 1513.8889
$ws.Range("J126").Value = 6721.1113
$ws.Range("K126").Value = 4541.6667
$ws.Range("L126").Value = 20163.3339
$ws.Range("M126").Value = -2071.6667
$ws.Range("N126").Value = -25103.3339
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

